$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.163.57'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.757.96'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.70'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.29'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.758.88'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.61%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +1.32%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.170'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.95%  '
$ws.Range("E11").Value = '  +1.83%  '
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("E13").Value = '  -0.44%  '
$ws.Range("E14").Value = '  +1.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.387.54'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.51%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.756.00'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.201.27'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.45%  '
$ws.Range("E18").Value = '  +1.47%  '
$ws.Range("E19").Value = '  -1.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.28'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.18'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +14.91%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '494.32'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.731'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.94%  '
$ws.Range("E24").Value = '  +7.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.98'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.32'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.36%  '
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("E30").Value = '  +1.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.19'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +3.55%  '
$ws.Range("E32").Value = '  +2.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.61'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.905.15'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.691.45'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.60%  '
$ws.Range("E36").Value = '  +0.13%  '
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("E38").Value = '  +1.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.99'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.66%  '
$ws.Range("E40").Value = '  +2.46%  '
$ws.Range("E41").Value = '  +0.99%  '
$ws.Range("E42").Value = '  +5.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '430.02'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.61%  '
$ws.Range("E44").Value = '  -0.71%  '
$ws.Range("E45").Value = '  +0.21%  '
$ws.Range("E46").Value = '  +1.11%  '
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.49'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.43'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.797.53'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.88%  '
$ws.Range("E51").Value = '  +0.57%  '
